# Apply the latest cryptos snapshot: updates the Price (D) and Volume(1h) (E) columns
# for rows 2-51. Numeric-looking Price values are written with a leading apostrophe so
# Excel keeps them as literal text (matching the workbook''s original formatting,
# e.g. "0.100" / "313.46") instead of auto-converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '40.816.85' }
    @{ Cell = 'E2'; Value = '  -2.26%  ' }
    @{ Cell = 'D3'; Value = '2.386.30' }
    @{ Cell = 'E3'; Value = '  -3.71%  ' }
    @{ Cell = 'E4'; Value = '  +0.10%  ' }
    @{ Cell = 'D5'; Value = '''313.46' }
    @{ Cell = 'E5'; Value = '  -1.87%  ' }
    @{ Cell = 'D6'; Value = '''88.14' }
    @{ Cell = 'E6'; Value = '  -5.68%  ' }
    @{ Cell = 'D7'; Value = '''0.530' }
    @{ Cell = 'E7'; Value = '  -4.31%  ' }
    @{ Cell = 'E8'; Value = '  +0.04%  ' }
    @{ Cell = 'D9'; Value = '''0.494' }
    @{ Cell = 'E9'; Value = '  -4.78%  ' }
    @{ Cell = 'D10'; Value = '''0.0826' }
    @{ Cell = 'E10'; Value = '  -5.05%  ' }
    @{ Cell = 'D11'; Value = '''31.23' }
    @{ Cell = 'E11'; Value = '  -6.28%  ' }
    @{ Cell = 'E12'; Value = '  -1.61%  ' }
    @{ Cell = 'D13'; Value = '2.754.86' }
    @{ Cell = 'E13'; Value = '  -3.72%  ' }
    @{ Cell = 'D14'; Value = '''6.62' }
    @{ Cell = 'E14'; Value = '  -4.20%  ' }
    @{ Cell = 'D15'; Value = '''15.22' }
    @{ Cell = 'E15'; Value = '  -3.98%  ' }
    @{ Cell = 'D16'; Value = '2.401.14' }
    @{ Cell = 'E16'; Value = '  -3.58%  ' }
    @{ Cell = 'D17'; Value = '''0.762' }
    @{ Cell = 'E17'; Value = '  -3.89%  ' }
    @{ Cell = 'D18'; Value = '40.737.66' }
    @{ Cell = 'E18'; Value = '  -2.36%  ' }
    @{ Cell = 'D19'; Value = '0.0₃0913' }
    @{ Cell = 'E19'; Value = '  -4.13%  ' }
    @{ Cell = 'D20'; Value = '''6.18' }
    @{ Cell = 'E20'; Value = '  -4.76%  ' }
    @{ Cell = 'D21'; Value = '''69.37' }
    @{ Cell = 'E21'; Value = '  -2.74%  ' }
    @{ Cell = 'D22'; Value = '''10.83' }
    @{ Cell = 'E22'; Value = '  -4.57%  ' }
    @{ Cell = 'D23'; Value = '''235.34' }
    @{ Cell = 'E23'; Value = '  -1.95%  ' }
    @{ Cell = 'D24'; Value = '''2.66' }
    @{ Cell = 'E24'; Value = '  -3.62%  ' }
    @{ Cell = 'E25'; Value = '  +0.18%  ' }
    @{ Cell = 'D26'; Value = '''1.83' }
    @{ Cell = 'E26'; Value = '  -5.90%  ' }
    @{ Cell = 'D27'; Value = '''23.69' }
    @{ Cell = 'E27'; Value = '  -4.40%  ' }
    @{ Cell = 'D28'; Value = '''2.21' }
    @{ Cell = 'E28'; Value = '  -2.54%  ' }
    @{ Cell = 'D29'; Value = '''9.42' }
    @{ Cell = 'E29'; Value = '  -4.25%  ' }
    @{ Cell = 'D30'; Value = '''34.08' }
    @{ Cell = 'E30'; Value = '  -5.94%  ' }
    @{ Cell = 'D31'; Value = '''155.50' }
    @{ Cell = 'E31'; Value = '  -1.80%  ' }
    @{ Cell = 'E32'; Value = '  +0.04%  ' }
    @{ Cell = 'D33'; Value = '''5.25' }
    @{ Cell = 'E33'; Value = '  -5.20%  ' }
    @{ Cell = 'D34'; Value = '''0.0733' }
    @{ Cell = 'E34'; Value = '  -4.54%  ' }
    @{ Cell = 'D35'; Value = '''2.42' }
    @{ Cell = 'E35'; Value = '  -6.72%  ' }
    @{ Cell = 'E36'; Value = '  -1.77%  ' }
    @{ Cell = 'D37'; Value = '''2.83' }
    @{ Cell = 'E37'; Value = '  -3.87%  ' }
    @{ Cell = 'D38'; Value = '''16.10' }
    @{ Cell = 'E38'; Value = '  -7.69%  ' }
    @{ Cell = 'D39'; Value = '''0.100' }
    @{ Cell = 'E39'; Value = '  -3.58%  ' }
    @{ Cell = 'D40'; Value = '''1.74' }
    @{ Cell = 'E40'; Value = '  -7.46%  ' }
    @{ Cell = 'D41'; Value = '''3.82' }
    @{ Cell = 'E41'; Value = '  -5.46%  ' }
    @{ Cell = 'E42'; Value = '  -7.85%  ' }
    @{ Cell = 'D43'; Value = '1.963.05' }
    @{ Cell = 'E43'; Value = '  -1.45%  ' }
    @{ Cell = 'D44'; Value = '''0.0272' }
    @{ Cell = 'E44'; Value = '  -4.92%  ' }
    @{ Cell = 'D45'; Value = '''17.66' }
    @{ Cell = 'E45'; Value = '  -7.31%  ' }
    @{ Cell = 'D46'; Value = '''2.81' }
    @{ Cell = 'E46'; Value = '  -6.46%  ' }
    @{ Cell = 'D47'; Value = '''9.38' }
    @{ Cell = 'E47'; Value = '  -0.78%  ' }
    @{ Cell = 'D48'; Value = '2.609.77' }
    @{ Cell = 'E48'; Value = '  -3.99%  ' }
    @{ Cell = 'D49'; Value = '''93.97' }
    @{ Cell = 'E49'; Value = '  -3.82%  ' }
    @{ Cell = 'D50'; Value = '''73.02' }
    @{ Cell = 'E50'; Value = '  -1.87%  ' }
    @{ Cell = 'D51'; Value = '''50.83' }
    @{ Cell = 'E51'; Value = '  -3.54%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
